$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("F2").Value = 3.9
$ws.Range("G2").Value = 6.4
$ws.Range("H2").Value = 1.69
$ws.Range("I2").Value = 1.88
$ws.Range("J2").Value = 3.75
$ws.Range("K2").Value = 4.9
$ws.Range("L2").Value = 1.01
$ws.Range("M2").Value = 1.01
$ws.Range("N2").Value = 2.1
$ws.Range("O2").Value = 1.23
$ws.Range("P2").Value = 2.1
$ws.Range("Q2").Value = 1.75
$ws.Range("R2").Value = 1.34
$ws.Range("S2").Value = 2.6
$ws.Range("T2").Value = 1.01
$ws.Range("U2").Value = 1.01
$ws.Range("V2").Value = 2.14
$ws.Range("W2").Value = 1.18
$ws.Range("X2").Value = 25
$ws.Range("Y2").Value = 13.5
$ws.Range("Z2").Value = 15.5
$ws.Range("AA2").Value = 25
$ws.Range("AB2").Value = 28
$ws.Range("AC2").Value = 13
$ws.Range("AD2").Value = 14.5
$ws.Range("AE2").Value = 25
$ws.Range("AF2").Value = 1000
$ws.Range("AG2").Value = 29
$ws.Range("AH2").Value = 26
$ws.Range("AI2").Value = 1000
$ws.Range("AJ2").Value = 1000
$ws.Range("AK2").Value = 1000
$ws.Range("AL2").Value = 1000
$ws.Range("AM2").Value = 1000
$ws.Range("AN2").Value = 1000
$ws.Range("AO2").Value = 1000
# Row 3
$ws.Range("F3").Value = 2.64
$ws.Range("G3").Value = 3.1
$ws.Range("H3").Value = 3.25
$ws.Range("I3").Value = 3.65
$ws.Range("J3").Value = 2.74
$ws.Range("K3").Value = 3.05
$ws.Range("P3").Value = 1.41
$ws.Range("Q3").Value = 3.1
# Row 4
$ws.Range("F4").Value = 1.65
$ws.Range("G4").Value = 1.8
$ws.Range("H4").Value = 5.2
$ws.Range("I4").Value = 6.4
$ws.Range("J4").Value = 3.75
$ws.Range("K4").Value = 4.5
$ws.Range("P4").Value = 1.97
$ws.Range("Q4").Value = 1.84
# Row 5
$ws.Range("F5").Value = 1.7
$ws.Range("G5").Value = 1.87
$ws.Range("H5").Value = 5.1
$ws.Range("I5").Value = 6.2
$ws.Range("J5").Value = 3.55
$ws.Range("K5").Value = 4.2
$ws.Range("M5").Value = 1.07
$ws.Range("P5").Value = 1.77
$ws.Range("Q5").Value = 2.1
$ws.Range("T5").Value = 1.97
$ws.Range("U5").Value = 1.84
$ws.Range("Y5").Value = 19
$ws.Range("AC5").Value = 9.800000000000001
$ws.Range("AF5").Value = 11.5
$ws.Range("AG5").Value = 970
$ws.Range("AI5").Value = 110
# Row 6
$ws.Range("F6").Value = 1.22
$ws.Range("G6").Value = 1.28
$ws.Range("H6").Value = 13.5
$ws.Range("I6").Value = 18.5
$ws.Range("J6").Value = 7
$ws.Range("K6").Value = 8.6
$ws.Range("P6").Value = 3
$ws.Range("Q6").Value = 1.39
# Row 7
$ws.Range("F7").Value = 5.6
$ws.Range("H7").Value = 1.17
$ws.Range("I7").Value = 1.21
$ws.Range("J7").Value = 8.800000000000001
$ws.Range("K7").Value = 11.5
$ws.Range("P7").Value = 3.1
$ws.Range("Q7").Value = 1.35
# Row 8
$ws.Range("F8").Value = 2.52
$ws.Range("G8").Value = 2.72
$ws.Range("H8").Value = 3.55
$ws.Range("I8").Value = 4.2
$ws.Range("J8").Value = 2.7
$ws.Range("K8").Value = 2.98
$ws.Range("P8").Value = 1.37
$ws.Range("Q8").Value = 3.25
# Row 9
$ws.Range("F9").Value = 1.73
$ws.Range("G9").Value = 1.85
$ws.Range("H9").Value = 5
$ws.Range("J9").Value = 3.35
$ws.Range("K9").Value = 4.1
$ws.Range("N9").Value = 2.78
$ws.Range("P9").Value = 1.61
$ws.Range("Q9").Value = 2.38
# Row 10
$ws.Range("F10").Value = 1.71
$ws.Range("G10").Value = 1.92
$ws.Range("H10").Value = 5.3
$ws.Range("I10").Value = 7
$ws.Range("J10").Value = 3.45
$ws.Range("K10").Value = 4.1
$ws.Range("P10").Value = 1.73
$ws.Range("Q10").Value = 2.1
# Row 11
$ws.Range("F11").Value = 2.78
$ws.Range("G11").Value = 4.4
$ws.Range("H11").Value = 2.26
$ws.Range("I11").Value = 2.82
$ws.Range("J11").Value = 2.74
$ws.Range("K11").Value = 3.5
$ws.Range("P11").Value = 1.54
$ws.Range("Q11").Value = 2.08
# Row 12
$ws.Range("F12").Value = 2.2
$ws.Range("G12").Value = 2.34
$ws.Range("H12").Value = 3.55
$ws.Range("I12").Value = 3.85
$ws.Range("J12").Value = 3.35
$ws.Range("K12").Value = 3.6
$ws.Range("P12").Value = 1.83
$ws.Range("Q12").Value = 2.06
# Row 13
$ws.Range("F13").Value = 2.16
$ws.Range("G13").Value = 2.68
$ws.Range("H13").Value = 2.8
$ws.Range("I13").Value = 4.6
$ws.Range("J13").Value = 2.84
$ws.Range("K13").Value = 3.7
$ws.Range("P13").Value = 1.73
$ws.Range("Q13").Value = 2.08
# Row 14
$ws.Range("F14").Value = 1.65
$ws.Range("G14").Value = 1.68
$ws.Range("H14").Value = 6.2
$ws.Range("I14").Value = 7
$ws.Range("J14").Value = 3.85
$ws.Range("P14").Value = 1.87
$ws.Range("Q14").Value = 2.02
